$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings are not
# auto-converted to numbers by Excel (preserves exact formatting, e.g. "1.00").
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "59.336.94"
$ws.Range("E2").Value = "  +2.12%  "

# Row 3
$ws.Range("D3").Value = "2.600.28"
$ws.Range("E3").Value = "  +1.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "535.06"
$ws.Range("E5").Value = "  +3.75%  "

# Row 6
$ws.Range("D6").Value = "140.56"
$ws.Range("E6").Value = "  +1.68%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  +1.75%  "

# Row 9
$ws.Range("D9").Value = "2.610.86"
$ws.Range("E9").Value = "  +1.37%  "

# Row 10
$ws.Range("E10").Value = "  -0.26%  "

# Row 11
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  +3.97%  "

# Row 12
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +3.53%  "

# Row 13
$ws.Range("D13").Value = "0.136"
$ws.Range("E13").Value = "  +2.65%  "

# Row 14
$ws.Range("D14").Value = "3.058.27"
$ws.Range("E14").Value = "  +1.37%  "

# Row 15
$ws.Range("D15").Value = "59.294.89"
$ws.Range("E15").Value = "  +2.17%  "

# Row 16
$ws.Range("D16").Value = "20.55"
$ws.Range("E16").Value = "  +1.93%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  +2.24%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.583.18"
$ws.Range("E18").Value = "  +0.67%  "

# Row 19
$ws.Range("D19").Value = "346.13"
$ws.Range("E19").Value = "  +3.75%  "

# Row 20
$ws.Range("D20").Value = "4.34"
$ws.Range("E20").Value = "  +1.49%  "

# Row 21
$ws.Range("D21").Value = "10.12"
$ws.Range("E21").Value = "  +0.74%  "

# Row 22
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +0.86%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
$ws.Range("D24").Value = "67.26"
$ws.Range("E24").Value = "  +2.20%  "

# Row 25
$ws.Range("E25").Value = "  +2.50%  "

# Row 26
$ws.Range("D26").Value = "0.408"
$ws.Range("E26").Value = "  +2.46%  "

# Row 27
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  +4.51%  "

# Row 29
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0738"
$ws.Range("E30").Value = "  +5.25%  "

# Row 31
$ws.Range("E31").Value = "  +5.44%  "

# Row 32
$ws.Range("D32").Value = "5.83"
$ws.Range("E32").Value = "  -0.70%  "

# Row 33
$ws.Range("D33").Value = "18.86"
$ws.Range("E33").Value = "  +1.58%  "

# Row 34
$ws.Range("D34").Value = "149.28"
$ws.Range("E34").Value = "  +0.33%  "

# Row 35
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  +3.38%  "

# Row 36
$ws.Range("E36").Value = "  +1.65%  "

# Row 37
$ws.Range("D37").Value = "36.94"
$ws.Range("E37").Value = "  +2.48%  "

# Row 38
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "0.845"
$ws.Range("E38").Value = "  +4.12%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  +5.05%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "0.838"
$ws.Range("E40").Value = "  +2.08%  "

# Row 41
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  +2.05%  "

# Row 42
$ws.Range("E42").Value = "  +0.18%  "

# Row 43
$ws.Range("D43").Value = "276.63"
$ws.Range("E43").Value = "  +2.50%  "

# Row 44
$ws.Range("D44").Value = "0.599"
$ws.Range("E44").Value = "  +2.52%  "

# Row 45
$ws.Range("D45").Value = "10.77"
$ws.Range("E45").Value = "  +0.43%  "

# Row 46
$ws.Range("D46").Value = "0.0962"
$ws.Range("E46").Value = "  +2.96%  "

# Row 47
$ws.Range("D47").Value = "0.0522"
$ws.Range("E47").Value = "  +2.16%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.946.87"
$ws.Range("E48").Value = "  -0.46%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0223"
$ws.Range("E49").Value = "  +2.86%  "

# Row 50
$ws.Range("D50").Value = "18.45"
$ws.Range("E50").Value = "  +5.97%  "

# Row 51
$ws.Range("D51").Value = "4.52"
$ws.Range("E51").Value = "  +2.46%  "
